$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.957.19'
$ws.Range("E2").Value = '  -1.65%  '
$ws.Range("D3").Value = '2.241.92'
$ws.Range("E3").Value = '  -1.84%  '
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").Value = '''246.75'
$ws.Range("E5").Value = '  -2.04%  '
$ws.Range("D6").Value = '''0.631'
$ws.Range("E6").Value = '  +0.60%  '
$ws.Range("D7").Value = '''75.31'
$ws.Range("E7").Value = '  +2.11%  '
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("D9").Value = '''0.631'
$ws.Range("E9").Value = '  -3.01%  '
$ws.Range("D10").Value = '''40.44'
$ws.Range("E10").Value = '  +2.90%  '
$ws.Range("D11").Value = '''0.0950'
$ws.Range("E11").Value = '  -3.28%  '
$ws.Range("D12").Value = '''7.18'
$ws.Range("E12").Value = '  -2.56%  '
$ws.Range("D13").Value = '''0.104'
$ws.Range("E13").Value = '  -2.50%  '
$ws.Range("D14").Value = '2.580.49'
$ws.Range("E14").Value = '  -1.42%  '
$ws.Range("D15").Value = '''14.84'
$ws.Range("E15").Value = '  -1.45%  '
$ws.Range("D16").Value = '''0.861'
$ws.Range("E16").Value = '  -2.10%  '
$ws.Range("D17").Value = '2.242.91'
$ws.Range("E17").Value = '  -2.11%  '
$ws.Range("D18").Value = '41.983.32'
$ws.Range("E18").Value = '  -1.32%  '
$ws.Range("D19").Value = '0.0₃0980'
$ws.Range("E19").Value = '  -2.82%  '
$ws.Range("D20").Value = '''6.13'
$ws.Range("E20").Value = '  -2.86%  '
$ws.Range("D21").Value = '''71.48'
$ws.Range("E21").Value = '  -1.07%  '
$ws.Range("D22").Value = '''2.25'
$ws.Range("E22").Value = '  +2.42%  '
$ws.Range("D23").Value = '''230.82'
$ws.Range("E23").Value = '  -0.97%  '
$ws.Range("E24").Value = '  +0.11%  '
$ws.Range("D25").Value = '''11.32'
$ws.Range("E25").Value = '  -1.88%  '
$ws.Range("E26").Value = '  -4.94%  '
$ws.Range("D27").Value = '''2.31'
$ws.Range("E27").Value = '  -4.32%  '
$ws.Range("D28").Value = '''7.27'
$ws.Range("E28").Value = '  +14.37%  '
$ws.Range("E29").Value = '  -1.54%  '
$ws.Range("D30").Value = '''169.57'
$ws.Range("E30").Value = '  +1.67%  '
$ws.Range("D31").Value = '''20.56'
$ws.Range("E31").Value = '  -2.55%  '
$ws.Range("D32").Value = '''34.15'
$ws.Range("E32").Value = '  +7.20%  '
$ws.Range("D33").Value = '''0.0842'
$ws.Range("E33").Value = '  +3.66%  '
$ws.Range("D34").Value = '''0.121'
$ws.Range("E34").Value = '  -5.52%  '
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("D36").Value = '''4.68'
$ws.Range("E36").Value = '  -0.43%  '
$ws.Range("D37").Value = '''4.91'
$ws.Range("E37").Value = '  +3.34%  '
$ws.Range("D38").Value = '''0.0300'
$ws.Range("E38").Value = '  -2.28%  '
$ws.Range("D39").Value = '''13.56'
$ws.Range("E39").Value = '  -2.55%  '
$ws.Range("D40").Value = '''5.94'
$ws.Range("E40").Value = '  +0.27%  '
$ws.Range("E41").Value = '  -6.64%  '
$ws.Range("D42").Value = '''110.23'
$ws.Range("E42").Value = '  +12.59%  '
$ws.Range("E43").Value = '  -4.58%  '
$ws.Range("D44").Value = '''60.41'
$ws.Range("E44").Value = '  -2.35%  '
$ws.Range("E45").Value = '  -3.96%  '
$ws.Range("D46").Value = '''0.100'
$ws.Range("E46").Value = '  -3.54%  '
$ws.Range("E47").Value = '  -0.46%  '
$ws.Range("D48").Value = '''1.13'
$ws.Range("E48").Value = '  -3.94%  '
$ws.Range("E49").Value = '  -1.18%  '
$ws.Range("E50").Value = '  -11.87%  '
$ws.Range("D51").Value = '''2.27'
$ws.Range("E51").Value = '  -1.34%  '
